{"js": "// 1) Update the build/revision date: two separate runs (\"2021-1\" + \"1-01\",\n//    together reading \"2021-11-01\") collapse into a single run with the\n//    corrected date \"2021-11-02\".\nconst dateResults = context.document.body.search(\"2021-11-01\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\"2021-11-02\", Word.InsertLocation.replace);\n}\n\n// 2) \"subspecialty\" -> \"sub-specialty\"\nconst subspecialtyResults = context.document.body.search(\"the subspecialty of cybersecurity\", { matchCase: true });\nsubspecialtyResults.load(\"items\");\nawait context.sync();\nif (subspecialtyResults.items.length > 0) {\n  subspecialtyResults.items[0].insertText(\"the sub-specialty of cybersecurity\", Word.InsertLocation.replace);\n}\n\n// 3) \", Zim wiki, MkDocs\" -> \", Zim personal wiki, MkDocs\"\nconst zimResults = context.document.body.search(\"Zim wiki, MkDocs\", { matchCase: true });\nzimResults.load(\"items\");\nawait context.sync();\nif (zimResults.items.length > 0) {\n  zimResults.items[0].insertText(\"Zim personal wiki, MkDocs\", Word.InsertLocation.replace);\n}\n\n// 4) Add \"Lua, Sql\" to the languages list, right after \"J (i.e. neo-APL),\"\nconst luaResults = context.document.body.search(\"J (i.e. neo-APL), Numpy\", { matchCase: true });\nluaResults.load(\"items\");\nawait context.sync();\nif (luaResults.items.length > 0) {\n  luaResults.items[0].insertText(\"J (i.e. neo-APL), Lua, Sql, Numpy\", Word.InsertLocation.replace);\n}\n\n// 5) Nudge the \"Computer Languages & Packages\" table's two column widths by\n//    2 twips: 1072 -> 1070 dxa (53.6 -> 53.5 pt) and 8899 -> 8901 dxa\n//    (444.95 -> 445.05 pt). Setting a cell's columnWidth resizes every row\n//    in that column plus the table's <w:gridCol> entry.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length > 0) {\n  const table = tables.items[0];\n  table.getCell(0, 0).columnWidth = 53.5; // 1070 dxa\n  table.getCell(0, 1).columnWidth = 445.05; // 8901 dxa\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Update the build/revision date: two separate runs (\"2021-1\" + \"1-01\",\n#    together reading \"2021-11-01\") collapse into a single run with the\n#    corrected date \"2021-11-02\".\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"2021-11-01\"\n$find.Replacement.Text = \"2021-11-02\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n\n# 2) \"subspecialty\" -> \"sub-specialty\"\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"the subspecialty of cybersecurity\"\n$find.Replacement.Text = \"the sub-specialty of cybersecurity\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n\n# 3) \", Zim wiki, MkDocs\" -> \", Zim personal wiki, MkDocs\"\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"Zim wiki, MkDocs\"\n$find.Replacement.Text = \"Zim personal wiki, MkDocs\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n\n# 4) Add \"Lua, Sql\" to the languages list, right after \"J (i.e. neo-APL),\"\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"J (i.e. neo-APL), Numpy\"\n$find.Replacement.Text = \"J (i.e. neo-APL), Lua, Sql, Numpy\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n\n# 5) Nudge the \"Computer Languages & Packages\" table's two column widths by\n#    2 twips: 1072 -> 1070 dxa (53.6 -> 53.5 pt) and 8899 -> 8901 dxa\n#    (444.95 -> 445.05 pt). Setting Column.Width resizes every row in that\n#    column plus the table's <w:gridCol> entry.\n$t = $d.Tables.Item(1)\n$t.Columns.Item(1).Width = 53.5\n$t.Columns.Item(2).Width = 445.05\n"}
